# Update the localization status text from "Ready for handoff" to
# "In Translation" everywhere it appears, then re-fit the Status columns
# so their width reflects the new (shorter) text, mirroring what Excel
# does automatically when a column is auto-sized after a content change.

$wb = $excel.ActiveWorkbook

$oldText = "Ready for handoff"
$newText = "In Translation"

foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace($oldText, $newText) | Out-Null
}

$overview = $wb.Worksheets.Item("Overview")
$overview.Columns("E:F").AutoFit() | Out-Null

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns("C:C").AutoFit() | Out-Null

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns("C:C").AutoFit() | Out-Null

Write-Host "Updated status text and auto-fit status columns."
